# The sheet holds one weekly price record per row for "Albahaca" at
# "Vega Modelo de Temuco". A new weekly record is inserted above the old
# row 263, pushing every following row down by one (old row 353 becomes
# row 354), and the sheet's used range grows from A1:R353 to A1:R354.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 263; rows 263-353 shift down to 264-354.
$ws.Rows(263).Insert()

# Populate the newly inserted row 263 with the new weekly record.
$ws.Range("A263").Value = 10
$ws.Range("B263").Value = "Vega Modelo de Temuco"
$ws.Range("C263").Value = "La Araucanía"
$ws.Range("D263").Value = 44985
$ws.Range("E263").Value = 9
$ws.Range("F263").Value = 100112052
$ws.Range("G263").Value = "Albahaca"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 40
$ws.Range("K263").Value = 6000
$ws.Range("L263").Value = 6000
$ws.Range("M263").Value = 6000
$ws.Range("N263").Value = "$/paquete"
$ws.Range("O263").Value = "Región del Maule"
$ws.Range("P263").Value = 6000
$ws.Range("Q263").Value = 1
$ws.Range("R263").Value = "Hortaliza"
